$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F5 cell value (last_update) from 1706219962 to 1706239962
$ws.Range("F5").Value = 1706239962

# Update the view: scroll so column C is the top-left visible column,
# and change the active selection to E6
$ws.Range("E6").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3 | Out-Null
